# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relabel countries whose shared-string slot shuffled around (the
#     underlying row data stayed put, only the name shown in column A
#     changed because of the reorder upstream) ---
$ws.Range("A53").Value = "Costa Rica"
$ws.Range("A54").Value = "Nepal"
$ws.Range("A55").Value = "Venezuela"

$ws.Range("A206").Value = "Timor Oriental"
$ws.Range("A207").Value = "Santa Lucia"

# --- Refresh the "last updated" timestamp caption ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Septiembre de 2020 a las 22:16"

# --- Updated per-country statistics (Casos totales, Nuevos casos,
#     Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Estados Unidos
$ws.Range("B4").Value = 7276752
$ws.Range("C4").Value = 32568
$ws.Range("D4").Value = 4505478
$ws.Range("E4").Value = 2562325
$ws.Range("G4").Value = 509
$ws.Range("H4").Value = 208949

# India
$ws.Range("B5").Value = 5990513
$ws.Range("C5").Value = 88942
$ws.Range("D5").Value = 4938641
$ws.Range("E5").Value = 957339
$ws.Range("G5").Value = 1123
$ws.Range("H5").Value = 94533

# Israel
$ws.Range("B27").Value = 227100
$ws.Range("C27").Value = 9201
$ws.Range("D27").Value = 157537
$ws.Range("E27").Value = 68122
$ws.Range("G27").Value = 29
$ws.Range("H27").Value = 1441

# Costa Rica (row 53)
$ws.Range("B53").Value = 72049
$ws.Range("C53").Value = 1233
$ws.Range("D53").Value = 27760
$ws.Range("E53").Value = 43461
$ws.Range("G53").Value = 16
$ws.Range("H53").Value = 828

# Nepal (row 54)
$ws.Range("B54").Value = 71821
$ws.Range("C54").Value = 1207
$ws.Range("D54").Value = 53013
$ws.Range("E54").Value = 18341
$ws.Range("G54").Value = 8
$ws.Range("H54").Value = 467

# Venezuela (row 55)
$ws.Range("B55").Value = 71273
$ws.Range("D55").Value = 60709
$ws.Range("E55").Value = 9973
$ws.Range("H55").Value = 591

# Costa de Marfil
$ws.Range("B85").Value = 19600
$ws.Range("C85").Value = 44
$ws.Range("D85").Value = 19122
$ws.Range("E85").Value = 358

# Malaui
$ws.Range("B116").Value = 5766
$ws.Range("C116").Value = 2
$ws.Range("D116").Value = 4185
$ws.Range("E116").Value = 1402

# Cabo Verde
$ws.Range("B117").Value = 5701
$ws.Range("C117").Value = 73
$ws.Range("D117").Value = 5018
$ws.Range("E117").Value = 627
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 56

# Siria
$ws.Range("B134").Value = 4038
$ws.Range("C134").Value = 37
$ws.Range("D134").Value = 1048
$ws.Range("E134").Value = 2802
$ws.Range("G134").Value = 3
$ws.Range("H134").Value = 188

# Yemen
$ws.Range("B153").Value = 2030
$ws.Range("C153").Value = 1
$ws.Range("D153").Value = 1260
$ws.Range("E153").Value = 183

# Togo
$ws.Range("B159").Value = 1736
$ws.Range("C159").Value = 14
$ws.Range("D159").Value = 1319
$ws.Range("E159").Value = 371
$ws.Range("G159").Value = 2
$ws.Range("H159").Value = 46

# Republica del Chad
$ws.Range("B167").Value = 1177
$ws.Range("C167").Value = 2
$ws.Range("D167").Value = 1005
$ws.Range("E167").Value = 89
